$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new entry to the list (A5) - new shared string with a trailing space
$ws.Range("A5").Value = "id osoby/potwora to name -4 "

# Move the selection to F13 (as recorded in the edited workbook)
$ws.Range("F13").Select()
